# Update cryptocurrency price (D) and volume-change (E) columns
# to reflect the latest scrape, per the GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number
# by Excel (e.g. "624.27") are temporarily switched to Text format so
# the value is stored as a string, matching the original cell type.
$textCells = @("D5","D6","D14","D18","D19","D21","D22","D25","D26","D34","D36","D41","D45","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price column (D)
$ws.Range("D2").Value = "69.384.28"
$ws.Range("D3").Value = "3.782.38"
$ws.Range("D5").Value = "624.27"
$ws.Range("D6").Value = "165.88"
$ws.Range("D7").Value = "3.779.89"
$ws.Range("D14").Value = "35.69"
$ws.Range("D15").Value = "4.418.22"
$ws.Range("D16").Value = "3.803.02"
$ws.Range("D17").Value = "69.381.04"
$ws.Range("D18").Value = "17.73"
$ws.Range("D19").Value = "7.13"
$ws.Range("D21").Value = "468.18"
$ws.Range("D22").Value = "9.62"
$ws.Range("D25").Value = "83.31"
$ws.Range("D26").Value = "12.04"
$ws.Range("D30").Value = "3.932.56"
$ws.Range("D34").Value = "28.83"
$ws.Range("D36").Value = "0.999"
$ws.Range("D37").Value = "3.734.11"
$ws.Range("D41").Value = "5.81"
$ws.Range("D45").Value = "43.35"
$ws.Range("D48").Value = "46.75"
$ws.Range("D49").Value = "151.98"

# Volume(1h) column (E)
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("E5").Value = "  +4.24%  "
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("E27").Value = "  +3.63%  "
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  +18.27%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("E40").Value = "  +7.63%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("E47").Value = "  +4.42%  "
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("E51").Value = "  +0.27%  "

# Restore the default (General) formatting on those cells so only the
# cell content changed, not its appearance/format.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
